# RPA datasets push 2023-12-10
# A new IPO row ("이닉스(구.이닉스정호)") is inserted at the top of the data
# table (row 2) and every subsequent row's data shifts down by one, which
# drops the former last row off the bottom of the table. Mechanically this
# is applied as a full rewrite of each data cell's value (rows 2-21),
# because the underlying XML keeps the same <row> elements/styles in place
# and only the cell contents change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    # Forces the cell to hold a literal text value even when the text looks
    # like a plain number (e.g. "19000"), matching the source data where
    # these "확정공모가" cells are stored as shared strings, not numbers.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2 - newly inserted IPO entry
$ws.Range("A2").Value = "이닉스(구.이닉스정호)"
$ws.Range("B2").Value = "2024.01.11~01.17"
$ws.Range("C2").Value = "9,200~11,000"
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = 27600
$ws.Range("F2").Value = "삼성증권"

# Row 3
$ws.Range("A3").Value = "에이치비인베스트먼트"
$ws.Range("B3").Value = "2024.01.08~01.12"
$ws.Range("C3").Value = "2,400~2,800"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = 16000
$ws.Range("F3").Value = "NH투자증권"

# Row 4
$ws.Range("A4").Value = "포스뱅크"
$ws.Range("B4").Value = "2024.01.05~01.11"
$ws.Range("C4").Value = "13,000~15,000"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = 19500
$ws.Range("F4").Value = "하나증권"

# Row 5
$ws.Range("A5").Value = "하나스팩30호"
$ws.Range("B5").Value = "2023.12.08~12.11"
$ws.Range("C5").Value = "2,000~2,000"
$ws.Range("D5").Value = "-"
$ws.Range("E5").Value = 14000
$ws.Range("F5").Value = "하나증권"

# Row 6
$ws.Range("A6").Value = "디에스단석(구,단석산업)"
$ws.Range("B6").Value = "2023.12.05~12.11"
$ws.Range("C6").Value = "79,000~89,000"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = 96380
$ws.Range("F6").Value = "KB증권,NH투자증권"

# Row 7
$ws.Range("A7").Value = "IBKS스팩23호"
$ws.Range("B7").Value = "2023.12.04~12.08"
$ws.Range("C7").Value = "2,000~2,000"
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = 8000
$ws.Range("F7").Value = "아이비케이투자증권"

# Row 8
$ws.Range("A8").Value = "블루엠텍"
$ws.Range("B8").Value = "2023.11.22~11.28"
$ws.Range("C8").Value = "15,000~19,000"
Set-TextValue "D8" "19000"
$ws.Range("E8").Value = 21000
$ws.Range("F8").Value = "하나증권,키움증권"

# Row 9
$ws.Range("A9").Value = "LS머트리얼즈"
$ws.Range("B9").Value = "2023.11.22~11.28"
$ws.Range("C9").Value = "4,400~5,500"
Set-TextValue "D9" "6000"
$ws.Range("E9").Value = 64350
$ws.Range("F9").Value = "키움증권,KB증권,이베스트투자증권,하이투자증권,NH투자증권"

# Row 10
$ws.Range("A10").Value = "삼성스팩9호"
$ws.Range("B10").Value = "2023.11.20~11.21"
$ws.Range("C10").Value = "2,000~2,000"
Set-TextValue "D10" "2000"
$ws.Range("E10").Value = 20000
$ws.Range("F10").Value = "삼성증권"

# Row 11
$ws.Range("A11").Value = "교보스팩15호"
$ws.Range("B11").Value = "2023.11.20~11.21"
$ws.Range("C11").Value = "2,000~2,000"
Set-TextValue "D11" "2000"
$ws.Range("E11").Value = 7000
$ws.Range("F11").Value = "교보증권"

# Row 12
$ws.Range("A12").Value = "케이엔에스"
$ws.Range("B12").Value = "2023.11.16~11.22"
$ws.Range("C12").Value = "19,000~22,000"
Set-TextValue "D12" "23000"
$ws.Range("E12").Value = 14250
$ws.Range("F12").Value = "신영증권"

# Row 13
$ws.Range("A13").Value = "NH스팩30호"
$ws.Range("B13").Value = "2023.11.15~11.16"
$ws.Range("C13").Value = "2,000~2,000"
Set-TextValue "D13" "2000"
$ws.Range("E13").Value = 16000
$ws.Range("F13").Value = "NH투자증권"

# Row 14
$ws.Range("A14").Value = "와이바이오로직스"
$ws.Range("B14").Value = "2023.11.10~11.16"
$ws.Range("C14").Value = "9,000~11,000"
Set-TextValue "D14" "9000"
$ws.Range("E14").Value = 13500
$ws.Range("F14").Value = "유안타증권"

# Row 15
$ws.Range("A15").Value = "에이텀"
$ws.Range("B15").Value = "2023.11.09~11.15"
$ws.Range("C15").Value = "23,000~30,000"
Set-TextValue "D15" "18000"
$ws.Range("E15").Value = 14950
$ws.Range("F15").Value = "하나증권"

# Row 16
$ws.Range("A16").Value = "에이에스텍"
$ws.Range("B16").Value = "2023.11.07~11.13"
$ws.Range("C16").Value = "21,000~25,000"
Set-TextValue "D16" "28000"
$ws.Range("E16").Value = 29547
$ws.Range("F16").Value = "미래에셋증권"

# Row 17
$ws.Range("A17").Value = "그린리소스"
$ws.Range("B17").Value = "2023.11.03~11.09"
$ws.Range("C17").Value = "11,000~14,000"
Set-TextValue "D17" "17000"
$ws.Range("E17").Value = 18040
$ws.Range("F17").Value = "NH투자증권"

# Row 18
$ws.Range("A18").Value = "한선엔지니어링"
$ws.Range("B18").Value = "2023.11.02~11.08"
$ws.Range("C18").Value = "5,200~6,000"
Set-TextValue "D18" "7000"
$ws.Range("E18").Value = 22100
$ws.Range("F18").Value = "대신증권"

# Row 19
$ws.Range("A19").Value = "에코아이"
$ws.Range("B19").Value = "2023.11.01~11.07"
$ws.Range("C19").Value = "28,500~34,700"
Set-TextValue "D19" "34700"
$ws.Range("E19").Value = 59251
$ws.Range("F19").Value = "KB증권"

# Row 20
$ws.Range("A20").Value = "동인기연(유가)"
$ws.Range("B20").Value = "2023.11.01~11.07"
$ws.Range("C20").Value = "33,000~37,000"
Set-TextValue "D20" "30000"
$ws.Range("E20").Value = 60654
$ws.Range("F20").Value = "NH투자증권"

# Row 21
$ws.Range("A21").Value = "스톰테크"
$ws.Range("B21").Value = "2023.10.31~11.06"
$ws.Range("C21").Value = "8,000~9,500"
Set-TextValue "D21" "11000"
$ws.Range("E21").Value = 26800
$ws.Range("F21").Value = "하이투자증권"
